# Add a new row for "Moaz Ashraf" with his email (as a mailto hyperlink) and
# his repo link, right under the existing header row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Moaz Ashraf"
$ws.Range("B2").Value = "moaazasm84@gmail.com"
$ws.Range("C2").Value = "https://github.com/moaaz311/Tradex.git"

# Turn the email cell into a mailto hyperlink (this also applies the
# built-in "Hyperlink" cell style to B2).
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:moaazasm84@gmail.com") | Out-Null

# Match the saved selection shown in the workbook after the edit.
$ws.Range("C2").Select() | Out-Null
